$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update measured odometry values (rotation test) that changed the deviation
# results below them.
$ws.Range("H15").Value = 63
$ws.Range("H16").Value = 59
$ws.Range("H17").Value = 91

# New block: consigne / Vcc / consigne cc computation
$ws.Range("G26").Value = "consigne rot gauche"
$ws.Range("H26").Value = 59.5

$ws.Range("G27").Value = "consigne rot droite"
$ws.Range("H27").Value = -60

$ws.Range("G28").Value = "Vcc"
$ws.Range("H28").Formula = '=(($H$26)-($H$27))/2'

$ws.Range("G29").Value = "Consigne gauche cc"
$ws.Range("H29").Formula = '=H28'

$ws.Range("G30").Value = "Consigne droite cc"
$ws.Range("H30").Formula = '=-H28'

# Selection / view state update (matches sheetView diff): scroll the window
# so row 11 is at the top and select H26, mirroring the author's on-screen
# state when they saved the file.
$ws.Activate()
$excel.Goto($ws.Range("A11"), $true)
$ws.Range("H26").Select()
